$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated naive component forecaster error-table values (rows 2-11, cols B-G)
$ws.Range("B2").Value = -0.100195032778358
$ws.Range("C2").Value = 1.879008151558315
$ws.Range("D2").Value = 15.52109306684352
$ws.Range("E2").Value = 3.939681848429327
$ws.Range("F2").Value = 4.026922141332304
$ws.Range("G2").Value = 23

$ws.Range("B3").Value = 0.1448133846575885
$ws.Range("C3").Value = 1.786984302705635
$ws.Range("D3").Value = 10.36886832397854
$ws.Range("E3").Value = 3.220072720293524
$ws.Range("F3").Value = 3.292514909155644
$ws.Range("G3").Value = 22

$ws.Range("B4").Value = -0.5106375609785742
$ws.Range("C4").Value = 1.177604939791601
$ws.Range("D4").Value = 4.283545554090664
$ws.Range("E4").Value = 2.069672813294571
$ws.Range("F4").Value = 2.055221296400692
$ws.Range("G4").Value = 21

$ws.Range("B5").Value = 0.07283602206394024
$ws.Range("C5").Value = 0.8438801149343232
$ws.Range("D5").Value = 1.589728322734128
$ws.Range("E5").Value = 1.260844289646476
$ws.Range("F5").Value = 1.291438706748583
$ws.Range("G5").Value = 20

$ws.Range("B6").Value = 0.009315935962682875
$ws.Range("C6").Value = 0.8519364303311395
$ws.Range("D6").Value = 1.743403430867843
$ws.Range("E6").Value = 1.320380032743544
$ws.Range("F6").Value = 1.356527761928607
$ws.Range("G6").Value = 19

$ws.Range("B7").Value = 0.01805680939183702
$ws.Range("C7").Value = 0.6865640439233759
$ws.Range("D7").Value = 0.7225809383130525
$ws.Range("E7").Value = 0.8500476094390552
$ws.Range("F7").Value = 0.8744944091926862
$ws.Range("G7").Value = 18

$ws.Range("B8").Value = 0.06286173468651945
$ws.Range("C8").Value = 0.8035601421598704
$ws.Range("D8").Value = 1.081400171720104
$ws.Range("E8").Value = 1.039903924273827
$ws.Range("F8").Value = 1.069948180945846
$ws.Range("G8").Value = 17

$ws.Range("B9").Value = 0.2080312578962194
$ws.Range("C9").Value = 0.6959573112162807
$ws.Range("D9").Value = 0.6635098639986297
$ws.Range("E9").Value = 0.8145611481028479
$ws.Range("F9").Value = 0.813376737467042
$ws.Range("G9").Value = 16

$ws.Range("B10").Value = 0.1458528612741254
$ws.Range("C10").Value = 0.6988170516431337
$ws.Range("D10").Value = 0.9094933485271182
$ws.Range("E10").Value = 0.9536736069154468
$ws.Range("F10").Value = 0.9755329814582326
$ws.Range("G10").Value = 15

$ws.Range("B11").Value = 0.2106889511877471
$ws.Range("C11").Value = 0.5147298635682541
$ws.Range("D11").Value = 0.3550355181284862
$ws.Range("E11").Value = 0.5958485697964595
$ws.Range("F11").Value = 0.578395630879239
$ws.Range("G11").Value = 14
